$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.147.89"
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").Value = "2.349.46"
$ws.Range("E3").Value = "  +6.76%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.69"
$ws.Range("E5").Value = "  +4.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.30"
$ws.Range("E6").Value = "  +3.36%  "
$ws.Range("E7").Value = "  +3.93%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.634"
$ws.Range("E9").Value = "  +7.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.41"
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0939"
$ws.Range("E11").Value = "  +4.33%  "
$ws.Range("E12").Value = "  +3.08%  "
$ws.Range("E13").Value = "  +13.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.105"
$ws.Range("E14").Value = "  +2.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.25"
$ws.Range("E15").Value = "  +9.89%  "
$ws.Range("D16").Value = "2.703.79"
$ws.Range("E16").Value = "  +6.70%  "
$ws.Range("D17").Value = "2.352.52"
$ws.Range("E17").Value = "  +5.85%  "
$ws.Range("D18").Value = "43.116.82"
$ws.Range("E18").Value = "  +3.02%  "
$ws.Range("E19").Value = "  +5.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.30"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.29"
$ws.Range("E21").Value = "  +4.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.58"
$ws.Range("E22").Value = "  +14.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.43"
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "253.33"
$ws.Range("E24").Value = "  +12.06%  "
$ws.Range("E25").Value = "  +2.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.06"
$ws.Range("E26").Value = "  +6.14%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.26"
$ws.Range("E28").Value = "  +4.42%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.53"
$ws.Range("E30").Value = "  +9.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.38"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.17"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0927"
$ws.Range("E33").Value = "  +7.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.98"
$ws.Range("E34").Value = "  +9.96%  "
$ws.Range("E35").Value = "  +7.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.01"
$ws.Range("E36").Value = "  +2.38%  "
$ws.Range("E37").Value = "  +5.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.12"
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("E39").Value = "  +3.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.71"
$ws.Range("E40").Value = "  +11.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.50"
$ws.Range("E41").Value = "  +3.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.48"
$ws.Range("E42").Value = "  +15.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.233"
$ws.Range("E43").Value = "  +3.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.79"
$ws.Range("E44").Value = "  +2.62%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.63"
$ws.Range("E46").Value = "  +5.40%  "
$ws.Range("E47").Value = "  +11.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.15"
$ws.Range("E48").Value = "  +8.67%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  +3.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.19"
$ws.Range("E51").Value = "  +6.61%  "
